$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for rows 2-50 ---
# NumberFormat is forced to Text ("@") before writing any D-column value that
# would otherwise be auto-parsed by Excel as a number (which would corrupt
# values like "0.0000278" into scientific notation, or turn "1.00" into 1).
$ws.Range("D2").Value = "67.790.14"
$ws.Range("E2").Value = "  -1.17%  "
$ws.Range("D3").Value = "3.785.96"
$ws.Range("E3").Value = "  -1.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.26"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.41"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "3.784.44"
$ws.Range("E7").Value = "  -1.80%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.52"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.454"
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000278"
$ws.Range("E13").Value = "  +3.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.62"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").Value = "4.419.77"
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("D16").Value = "3.754.82"
$ws.Range("E16").Value = "  -2.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.63"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").Value = "67.733.27"
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.20"
$ws.Range("E19").Value = "  -2.72%  "
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("E21").Value = "  -5.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "468.89"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.720"
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000148"
$ws.Range("E24").Value = "  -8.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.89"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.21"
$ws.Range("E26").Value = "  -1.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.20"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.28"
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.92"
$ws.Range("E30").Value = "  -1.93%  "
$ws.Range("D31").Value = "3.931.86"
$ws.Range("E31").Value = "  -1.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.65"
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.57"
$ws.Range("E33").Value = "  -3.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.23"
$ws.Range("E34").Value = "  -3.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.15"
$ws.Range("E35").Value = "  -2.27%  "
$ws.Range("D36").Value = "3.746.62"
$ws.Range("E36").Value = "  -1.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.82"
$ws.Range("E37").Value = "  +2.61%  "
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.01"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.81"
$ws.Range("E41").Value = "  -2.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.997"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.312"
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.95"
$ws.Range("E46").Value = "  -1.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.81"
$ws.Range("E47").Value = "  -2.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "396.72"
$ws.Range("E48").Value = "  -5.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000271"
$ws.Range("E49").Value = "  -8.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "140.69"
$ws.Range("E50").Value = "  -0.82%  "

# --- Row 51: coin swapped from VeChain to Arweave ---
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "39.24"
$ws.Range("E51").Value = "  +3.24%  "
